# Generate Report for Handoff
# Update the "bb9e966a-6ebc-4b9f-b344-25f12d04e8ac.md" row with the latest
# handoff timestamps on the Overview sheet and on each locale's detail sheet.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest Handoff Date" for the bb9e966a... row (row 5)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D5").Value = "2016-03-22 06:12:01"

# zh-cn detail sheet: "Latest Handoff Datetime" for the bb9e966a... row (row 5)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E5").Value = "2016-03-22 06:11:53"

# de-de detail sheet: "Latest Handoff Datetime" for the bb9e966a... row (row 5)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E5").Value = "2016-03-22 06:12:01"
